$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DetectiveWork")
$ws.Range("B2").Value = "CAFR2017 p63"
$ws.Range("B4").Value = "FY2016 end"
$ws.Range("B5").Value = "FY2017 end"
$ws.Range("C3").Value = "Pension Liability"
$ws.Range("D3").Value = "Fiduciary net position"
$ws.Range("F3").Value = "AA"
$ws.Range("E3").Value = "MA"
$ws.Range("G3").Value = "TDA assets"
$ws.Range("C4").Value = 70000777
$ws.Range("D4").Value = 43629545
$ws.Range("G4").Value = 20292733
$ws.Range("C5").Value = 73323430
$ws.Range("D5").Value = 50095723
$ws.Range("G5").Value = 22004183
